$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.959.10"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "2.931.87"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.568"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.634"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "3.395.88"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "2.928.42"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.989"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "52.004.40"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("E26").Value = "  +11.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.45%  "
$ws.Range("E30").Value = "  +14.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "52.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  -14.75%  "
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("E40").Value = "  -3.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("D48").Value = "2.143.34"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("E49").Value = "  -7.35%  "
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("E51").Value = "  +0.56%  "
